$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) contains values that look numeric (e.g. "15.03",
# "0.5290") as well as multi-dot strings (e.g. "27.236.97"). The source
# workbook stores all of these as plain text (inlineStr) cells. Temporarily
# force the Price column to a text number format so Excel's COM layer does
# not auto-convert the assigned strings into floating point numbers
# (which would also silently drop significant/trailing digits).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.236.97"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "1.820.56"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "313.23"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.4464"
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("D8").Value = "0.3772"
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("D9").Value = "0.07396"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").Value = "0.8787"
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("D11").Value = "20.84"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "1.821.01"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "6.710"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "5.415"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "93.09"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "0.07116"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "0.000008795"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "15.03"
$ws.Range("D21").Value = "27.236.19"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").Value = "5.354"
$ws.Range("E22").Value = "  +3.68%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "1.964"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").Value = "151.05"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "2.289"
$ws.Range("E26").Value = "  +2.92%  "
$ws.Range("D27").Value = "18.57"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "5.338"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").Value = "117.37"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").Value = "0.08867"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "0.7823"
$ws.Range("E31").Value = "  +4.17%  "
$ws.Range("D32").Value = "1.192"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "4.566"
$ws.Range("D34").Value = "2.907"
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("D35").Value = "1.001"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "1.106"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").Value = "0.01972"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Value = "0.05256"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").Value = "7.313"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("D40").Value = "0.5290"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").Value = "2.862"
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").Value = "0.1704"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "2.269"
$ws.Range("E43").Value = "  +14.69%  "
$ws.Range("D44").Value = "8.588"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").Value = "0.5027"
$ws.Range("E45").Value = "  -3.75%  "
$ws.Range("D46").Value = "10.58"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "104.84"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "1.685"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "0.06387"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "66.02"
$ws.Range("E51").Value = "  +4.82%  "

# Restore the original (default/"Normal") style on the Price column so the
# cells keep the same unstyled appearance they had before this script ran.
$priceRange.Style = "Normal"
